$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values only
$ws.Range("B3").Value = 1030526647407.139
$ws.Range("C3").Value = 1028875162395.27
$ws.Range("D3").Value = 1635220702709212

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 1080346445265.233
$ws.Range("C4").Value = 970247444473.8754
$ws.Range("D4").Value = 212894459750680.4

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 108832131653093.8
$ws.Range("C5").Value = 57837435133899.56
$ws.Range("D5").Value = 479098480933519.8
